{"js": "// Replace each old two-digit multiplication expression with its new value.\n// The source document contains a grid table of \"AA\u00d7BB=\" cells; every old\n// value below is unique in the document, so a direct search + replace\n// for each pair reproduces the diff exactly.\nconst replacements = [\n  [\"62\u00d758=\", \"55\u00d766=\"],\n  [\"22\u00d739=\", \"95\u00d786=\"],\n  [\"79\u00d773=\", \"68\u00d793=\"],\n  [\"39\u00d724=\", \"54\u00d751=\"],\n  [\"74\u00d734=\", \"30\u00d732=\"],\n  [\"43\u00d721=\", \"21\u00d749=\"],\n  [\"45\u00d716=\", \"81\u00d766=\"],\n  [\"54\u00d771=\", \"72\u00d794=\"],\n  [\"40\u00d713=\", \"27\u00d738=\"],\n  [\"35\u00d733=\", \"94\u00d755=\"],\n  [\"94\u00d753=\", \"82\u00d711=\"],\n  [\"65\u00d725=\", \"31\u00d736=\"],\n  [\"50\u00d790=\", \"58\u00d774=\"],\n  [\"17\u00d765=\", \"56\u00d774=\"],\n  [\"31\u00d725=\", \"41\u00d733=\"],\n  [\"18\u00d753=\", \"29\u00d783=\"],\n  [\"59\u00d780=\", \"64\u00d794=\"],\n  [\"46\u00d732=\", \"65\u00d740=\"],\n  [\"21\u00d727=\", \"44\u00d756=\"],\n  [\"71\u00d793=\", \"42\u00d768=\"],\n  [\"64\u00d772=\", \"66\u00d798=\"],\n  [\"44\u00d723=\", \"58\u00d732=\"],\n  [\"25\u00d799=\", \"77\u00d773=\"],\n  [\"81\u00d792=\", \"76\u00d793=\"],\n  [\"84\u00d723=\", \"69\u00d727=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old two-digit multiplication expression with its new value.\n# The source document contains a grid table of \"AA\u00d7BB=\" cells; every old\n# value below is unique in the document, so a Find/Replace pass for each\n# pair reproduces the diff exactly. $replacements is a flat list of\n# old/new string pairs (even index = search text, odd index = replacement).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    \"62\u00d758=\", \"55\u00d766=\",\n    \"22\u00d739=\", \"95\u00d786=\",\n    \"79\u00d773=\", \"68\u00d793=\",\n    \"39\u00d724=\", \"54\u00d751=\",\n    \"74\u00d734=\", \"30\u00d732=\",\n    \"43\u00d721=\", \"21\u00d749=\",\n    \"45\u00d716=\", \"81\u00d766=\",\n    \"54\u00d771=\", \"72\u00d794=\",\n    \"40\u00d713=\", \"27\u00d738=\",\n    \"35\u00d733=\", \"94\u00d755=\",\n    \"94\u00d753=\", \"82\u00d711=\",\n    \"65\u00d725=\", \"31\u00d736=\",\n    \"50\u00d790=\", \"58\u00d774=\",\n    \"17\u00d765=\", \"56\u00d774=\",\n    \"31\u00d725=\", \"41\u00d733=\",\n    \"18\u00d753=\", \"29\u00d783=\",\n    \"59\u00d780=\", \"64\u00d794=\",\n    \"46\u00d732=\", \"65\u00d740=\",\n    \"21\u00d727=\", \"44\u00d756=\",\n    \"71\u00d793=\", \"42\u00d768=\",\n    \"64\u00d772=\", \"66\u00d798=\",\n    \"44\u00d723=\", \"58\u00d732=\",\n    \"25\u00d799=\", \"77\u00d773=\",\n    \"81\u00d792=\", \"76\u00d793=\",\n    \"84\u00d723=\", \"69\u00d727=\"\n)\n\nfor ($i = 0; $i -lt $replacements.Count; $i += 2) {\n    $oldText = $replacements[$i]\n    $newText = $replacements[$i + 1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
